$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32-43 were reordered upstream (each row now carries the field values
# that used to belong to a different row in this block). Re-apply that new
# order by writing out the permuted values for every affected column
# (A, B, D, E, F, G, H, I, J, Q, R); all other columns are left untouched.

# Row 32 now holds the record previously stored in row 41
$ws.Range("A32").Value = 112017392
$ws.Range("B32").Value = 90858
$ws.Range("E32").Value = 5449
$ws.Range("Q32").Value = 682712
$ws.Range("R32").Value = 6575458
$ws.Range("D32").Value = "NT"
$ws.Range("F32").Value = "Svart taggsvamp"
$ws.Range("G32").Value = "Phellodon niger"
$ws.Range("H32").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("J32").Value = ""
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = ""

# Row 33 now holds the record previously stored in row 39
$ws.Range("A33").Value = 112017512
$ws.Range("B33").Value = 88180
$ws.Range("E33").Value = 6276
$ws.Range("Q33").Value = 683037
$ws.Range("R33").Value = 6575484
$ws.Range("D33").Value = "VU"
$ws.Range("F33").Value = "Goliatmusseron"
$ws.Range("G33").Value = "Tricholoma matsutake"
$ws.Range("H33").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("J33").Value = "fruktkroppar"
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = "4"

# Row 34 now holds the record previously stored in row 38
$ws.Range("A34").Value = 112017413
$ws.Range("B34").Value = 90857
$ws.Range("E34").Value = 5448
$ws.Range("Q34").Value = 682734
$ws.Range("R34").Value = 6575482
$ws.Range("D34").Value = "NT"
$ws.Range("F34").Value = "Svartvit taggsvamp"
$ws.Range("G34").Value = "Phellodon connatus"
$ws.Range("H34").Value = "(Schultz) nom.prov"
$ws.Range("J34").Value = ""
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value = ""

# Row 35 now holds the record previously stored in row 40
$ws.Range("A35").Value = 112017488
$ws.Range("B35").Value = 90826
$ws.Range("E35").Value = 4366
$ws.Range("Q35").Value = 682956
$ws.Range("R35").Value = 6575474
$ws.Range("D35").Value = "LC"
$ws.Range("F35").Value = "Skarp dropptaggsvamp"
$ws.Range("G35").Value = "Hydnellum peckii"
$ws.Range("H35").Value = "Banker"
$ws.Range("J35").Value = ""
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = ""

# Row 36 now holds the record previously stored in row 43
$ws.Range("A36").Value = 112017130
$ws.Range("B36").Value = 90814
$ws.Range("E36").Value = 4364
$ws.Range("Q36").Value = 682695
$ws.Range("R36").Value = 6575454
$ws.Range("D36").Value = "LC"
$ws.Range("F36").Value = "Dropptaggsvamp"
$ws.Range("G36").Value = "Hydnellum ferrugineum"
$ws.Range("H36").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("J36").Value = ""
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = ""

# Row 37 now holds the record previously stored in row 35
$ws.Range("A37").Value = 112017326
$ws.Range("B37").Value = 90808
$ws.Range("E37").Value = 4362
$ws.Range("Q37").Value = 682714
$ws.Range("R37").Value = 6575496
$ws.Range("D37").Value = "NT"
$ws.Range("F37").Value = "Blå taggsvamp"
$ws.Range("G37").Value = "Hydnellum caeruleum"
$ws.Range("H37").Value = "(Hornem.) P.Karst."
$ws.Range("J37").Value = ""
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = ""

# Row 38 now holds the record previously stored in row 37
$ws.Range("A38").Value = 112017159
$ws.Range("B38").Value = 90858
$ws.Range("E38").Value = 5449
$ws.Range("Q38").Value = 682699
$ws.Range("R38").Value = 6575482
$ws.Range("D38").Value = "NT"
$ws.Range("F38").Value = "Svart taggsvamp"
$ws.Range("G38").Value = "Phellodon niger"
$ws.Range("H38").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("J38").Value = ""
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = ""

# Row 39 now holds the record previously stored in row 36
$ws.Range("A39").Value = 112017252
$ws.Range("B39").Value = 90814
$ws.Range("E39").Value = 4364
$ws.Range("Q39").Value = 682711
$ws.Range("R39").Value = 6575494
$ws.Range("D39").Value = "LC"
$ws.Range("F39").Value = "Dropptaggsvamp"
$ws.Range("G39").Value = "Hydnellum ferrugineum"
$ws.Range("H39").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("J39").Value = ""
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = ""

# Row 40 now holds the record previously stored in row 42
$ws.Range("A40").Value = 112017224
$ws.Range("B40").Value = 90826
$ws.Range("E40").Value = 4366
$ws.Range("Q40").Value = 682703
$ws.Range("R40").Value = 6575491
$ws.Range("D40").Value = "LC"
$ws.Range("F40").Value = "Skarp dropptaggsvamp"
$ws.Range("G40").Value = "Hydnellum peckii"
$ws.Range("H40").Value = "Banker"
$ws.Range("J40").Value = ""
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = ""

# Row 41 now holds the record previously stored in row 32
$ws.Range("A41").Value = 112017534
$ws.Range("B41").Value = 88140
$ws.Range("E41").Value = 1593
$ws.Range("Q41").Value = 683073
$ws.Range("R41").Value = 6575478
$ws.Range("D41").Value = "VU"
$ws.Range("F41").Value = "Lakritsmusseron"
$ws.Range("G41").Value = "Tricholoma apium"
$ws.Range("H41").Value = "Jul.Schäff."
$ws.Range("J41").Value = "fruktkroppar"
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = "4"

# Row 42 now holds the record previously stored in row 33
$ws.Range("A42").Value = 112017430
$ws.Range("B42").Value = 90857
$ws.Range("E42").Value = 5448
$ws.Range("Q42").Value = 682793
$ws.Range("R42").Value = 6575520
$ws.Range("D42").Value = "NT"
$ws.Range("F42").Value = "Svartvit taggsvamp"
$ws.Range("G42").Value = "Phellodon connatus"
$ws.Range("H42").Value = "(Schultz) nom.prov"
$ws.Range("J42").Value = ""
$ws.Range("I42").NumberFormat = "@"
$ws.Range("I42").Value = ""

# Row 43 now holds the record previously stored in row 34
$ws.Range("A43").Value = 112017465
$ws.Range("B43").Value = 88180
$ws.Range("E43").Value = 6276
$ws.Range("Q43").Value = 682896
$ws.Range("R43").Value = 6575514
$ws.Range("D43").Value = "VU"
$ws.Range("F43").Value = "Goliatmusseron"
$ws.Range("G43").Value = "Tricholoma matsutake"
$ws.Range("H43").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("J43").Value = "fruktkroppar"
$ws.Range("I43").NumberFormat = "@"
$ws.Range("I43").Value = "3"
